$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the five main-dish entries in column A (rows 2-6) -------------
# Old dishes (Pata y muslo con ensalada, Pechuga grille gratinada con pure
# de papas, Canelones mixtos con crema de espinaca, Sandwich de milanesa
# con lechuga y tomate, Sorrentinos con crema de Roquefort y jamon) are
# replaced by the new ACT menu entries. Replacing the cell values causes
# the old, now-unreferenced shared strings to be dropped automatically and
# the new ones to be appended, which also re-numbers the (unchanged)
# salad-ingredient strings in column H.
$ws.Range("A2").Value = "Filet de Merluza con pure mixto"
$ws.Range("A3").Value = "Guiso de lentejas"
$ws.Range("A4").Value = "Sandwich primavera en pan multicereal`n"
$ws.Range("A5").Value = "Ravioles Piamontés con salsa 4 quesos"
$ws.Range("A6").Value = "Sorrentinos veganos de berenjenas y quinoa con salsa rosa"

# --- Re-style the new A2 text (it was pasted in from elsewhere, carrying
# its own small, dark-grey Arial font). Build the font on a scratch cell
# that starts from the default (unstyled) format and copy just the
# formatting over, so A2 ends up with a plain "apply font only" style
# instead of inheriting the wrap/vertical-center alignment of style 2. ---
$helper = $ws.Range("Z1")
$helper.Font.Name = "Arial"
$helper.Font.Size = 9
$helper.Font.Color = 2236962
$helper.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$helper.Clear()
$excel.CutCopyMode = 0

# --- Row heights adjust to fit the new text --------------------------
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 72
$ws.Rows.Item(6).RowHeight = 72

# --- Update the sheet view: scroll back to the top and move the
# selection to K5 -------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$ws.Range("K5").Select()

# --- Page setup: portrait, paper size 9 (A4) ---------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
